# Refresh the crypto price/volume snapshot (symbol list update, 2022-12-25).
#
# All data cells in this sheet are stored as text (inlineStr), including the
# numeric-looking "Price" column D. If we just assign a numeric-looking
# string to .Value, Excel will happily reinterpret it as a number and drop
# meaningful trailing/leading zeros (e.g. "243.00" -> 243, "0.006170" ->
# 0.00617). To keep the text formatting exactly as scraped, force the
# NumberFormat to Text ("@") on those cells before writing the new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Price (column D) refreshes for rows whose coin/rank didn't move ---
Set-TextValue "D2"  "243.00"
Set-TextValue "D3"  "23.09"
Set-TextValue "D4"  "5.391"
Set-TextValue "D6"  "3.403"
Set-TextValue "D7"  "6.481"
Set-TextValue "D8"  "0.8089"
Set-TextValue "D9"  "0.9082"
Set-TextValue "D10" "0.1418"
Set-TextValue "D11" "0.07428"
Set-TextValue "D12" "0.03313"
Set-TextValue "D13" "0.03069"
Set-TextValue "D14" "0.09339"
Set-TextValue "D15" "3.852"
Set-TextValue "D16" "0.001593"
Set-TextValue "D17" "0.04633"
Set-TextValue "D18" "0.0005937"
Set-TextValue "D19" "0.006120"
Set-TextValue "D20" "0.005033"

# --- Rows 21-23: ranking reshuffled (UpBots/BitKan/NitroEx rotate) ---
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D21" "0.0009862"
$ws.Range("E21").Value = "20BitKanKAN"

$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D22" "0.00007797"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").Value = "UpBots"
$ws.Range("C23").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue "D23" "0.0002899"
$ws.Range("E23").Value = "22UpBotsUBXT"

# --- Remaining Price (and one label) refreshes further down the sheet ---
Set-TextValue "D26" "0.3215"

Set-TextValue "D40" "0.03891"

Set-TextValue "D41" "0.006170"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

Set-TextValue "D42" "0.1071"
Set-TextValue "D43" "0.002799"
Set-TextValue "D44" "0.007162"
Set-TextValue "D45" "0.00005188"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.0005797"
Set-TextValue "D48" "1.045"
Set-TextValue "D49" "0.002263"
Set-TextValue "D50" "0.00002099"
Set-TextValue "D51" "0.0001999"
